$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.797.94"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "3.756.57"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'613.77"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'178.92"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("D7").Value = "3.766.65"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  -2.00%  "
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("D11").Value = "'6.60"
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("E12").Value = "  -3.17%  "
$ws.Range("D13").Value = "'40.19"
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "4.378.20"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "3.750.29"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "69.857.10"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("D19").Value = "'7.47"
$ws.Range("E19").Value = "  -1.53%  "
$ws.Range("D20").Value = "'16.46"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("D21").Value = "'502.80"
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("D22").Value = "'9.22"
$ws.Range("E22").Value = "  -3.68%  "
$ws.Range("D23").Value = "'0.723"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").Value = "'2.65"
$ws.Range("E24").Value = "  +6.16%  "
$ws.Range("D25").Value = "'86.10"
$ws.Range("E25").Value = "  -2.26%  "
$ws.Range("D26").Value = "'11.50"
$ws.Range("E26").Value = "  +4.09%  "
$ws.Range("D27").Value = "'12.97"
$ws.Range("E27").Value = "  -3.53%  "
$ws.Range("E28").Value = "  +7.00%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'8.21"
$ws.Range("E30").Value = "  +4.23%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'2.49"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("E32").Value = "  +2.43%  "
$ws.Range("D33").Value = "'30.51"
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("D38").Value = "'0.354"
$ws.Range("E38").Value = "  +4.56%  "
$ws.Range("E39").Value = "  +4.01%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'452.69"
$ws.Range("E40").Value = "  +7.12%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'3.08"
$ws.Range("E41").Value = "  +13.32%  "
$ws.Range("E42").Value = "  -4.78%  "
$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D43").Value = "'45.73"
$ws.Range("E43").Value = "  +3.46%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'49.77"
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("D45").Value = "'8.59"
$ws.Range("E45").Value = "  -2.53%  "
$ws.Range("D46").Value = "2.956.19"
$ws.Range("E46").Value = "  -4.48%  "
$ws.Range("D47").Value = "'0.0361"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("D48").Value = "'27.21"
$ws.Range("E48").Value = "  -2.43%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'138.85"
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").Value = "'2.50"
$ws.Range("E51").Value = "  -1.23%  "
